# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) updates: row -> new value for column F
$sheet1Updates = @{
    3  = 157
    5  = 518
    6  = 1556
    7  = 5
    8  = 1192
    9  = 116
    10 = 219
    11 = 160
    13 = 3
    14 = 4
    15 = 230
    16 = 128
    17 = 192
    18 = 177
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (sheet4) updates: row -> new value for column F
$sheet4Updates = @{
    3  = 157
    5  = 518
    6  = 1556
    8  = 5
    9  = 1192
    10 = 116
    11 = 219
    12 = 160
    14 = 3
    15 = 4
    16 = 230
    17 = 128
    18 = 192
    19 = 177
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
